$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 18.13101768493652
$ws.Range("C3").Value = 16.71600341796875
$ws.Range("C4").Value = 16.2808895111084
$ws.Range("C5").Value = 16.20316505432129
$ws.Range("C6").Value = 16.35098457336426
